$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 923.5  # Dont Be So Tallow
$ws.Range("I12").Value = 864.1429000000001  # Dont Be So Tallow
$ws.Range("K12").Value = 864.1429000000001  # Dont Be So Tallow
$ws.Range("M12").Value = -694.1429000000001  # Dont Be So Tallow
$ws.Range("H18").Value = 15040.75  # You Grow, Girl
$ws.Range("I18").Value = 9833.333000000001  # You Grow, Girl
$ws.Range("K18").Value = 9833.333000000001  # You Grow, Girl
$ws.Range("M18").Value = -9549.333000000001  # You Grow, Girl
$ws.Range("H33").Value = 71.916664  # Glazed and Confused
$ws.Range("I33").Value = 54.22222  # Glazed and Confused
$ws.Range("K33").Value = 54.22222  # Glazed and Confused
$ws.Range("M33").Value = 174.77778  # Glazed and Confused
$ws.Range("H53").Value = 7812.6924  # No Accounting for Waste
$ws.Range("I53").Value = 12518.625  # No Accounting for Waste
$ws.Range("J53").Value = 283.2  # No Accounting for Waste
$ws.Range("K53").Value = 12518.625  # No Accounting for Waste
$ws.Range("L53").Value = 283.2  # No Accounting for Waste
$ws.Range("M53").Value = -11881.625  # No Accounting for Waste
$ws.Range("N53").Value = -1557.2  # No Accounting for Waste
$ws.Range("H112").Value = 2747.4167  # Making Ends Meet
$ws.Range("I112").Value = 400  # Making Ends Meet
$ws.Range("K112").Value = 1200  # Making Ends Meet
$ws.Range("M112").Value = -92  # Making Ends Meet
$ws.Range("H132").Value = 1197.4857  # Fast-forwarding Flora
$ws.Range("I132").Value = 1144.4706  # Fast-forwarding Flora
$ws.Range("K132").Value = 3433.4118  # Fast-forwarding Flora
$ws.Range("M132").Value = -903.4118000000003  # Fast-forwarding Flora
$ws.Range("H138").Value = 3153.1924  # All-night Crafting
$ws.Range("I138").Value = 4186  # All-night Crafting
$ws.Range("J138").Value = 2120.3845  # All-night Crafting
$ws.Range("K138").Value = 12558  # All-night Crafting
$ws.Range("L138").Value = 6361.1535  # All-night Crafting
$ws.Range("M138").Value = -7418  # All-night Crafting
$ws.Range("N138").Value = -16641.1535  # All-night Crafting

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2734.4417  # Ingot We Trust
$ws.Range("I32").Value = 1607.2878  # Ingot We Trust
$ws.Range("K32").Value = 1607.2878  # Ingot We Trust
$ws.Range("M32").Value = -1320.2878  # Ingot We Trust
$ws.Range("H122").Value = 1652.9667  # Haste for High Durium
$ws.Range("I122").Value = 1628.56  # Haste for High Durium
$ws.Range("J122").Value = 1775  # Haste for High Durium
$ws.Range("K122").Value = 4885.68  # Haste for High Durium
$ws.Range("L122").Value = 5325  # Haste for High Durium
$ws.Range("M122").Value = -2435.68  # Haste for High Durium
$ws.Range("N122").Value = -10225  # Haste for High Durium

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 999.5  # High Steal
$ws.Range("I94").Value = 999.5  # High Steal
$ws.Range("K94").Value = 999.5  # High Steal
$ws.Range("M94").Value = -548.5  # High Steal
$ws.Range("H107").Value = 2100.8462  # The Gold Experience
$ws.Range("I107").Value = 1647.2222  # The Gold Experience
$ws.Range("J107").Value = 3121.5  # The Gold Experience
$ws.Range("K107").Value = 1647.2222  # The Gold Experience
$ws.Range("L107").Value = 3121.5  # The Gold Experience
$ws.Range("M107").Value = 272.7778000000001  # The Gold Experience
$ws.Range("N107").Value = -6961.5  # The Gold Experience

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1157.8  # You Do the Heavy Lifting
$ws.Range("I58").Value = 989.2  # You Do the Heavy Lifting
$ws.Range("J58").Value = 1326.4  # You Do the Heavy Lifting
$ws.Range("K58").Value = 989.2  # You Do the Heavy Lifting
$ws.Range("L58").Value = 1326.4  # You Do the Heavy Lifting
$ws.Range("M58").Value = -786.2  # You Do the Heavy Lifting
$ws.Range("N58").Value = -1732.4  # You Do the Heavy Lifting
$ws.Range("H132").Value = 1923.3914  # Hull Lotta Damage
$ws.Range("I132").Value = 1064.5264  # Hull Lotta Damage
$ws.Range("K132").Value = 3193.5792  # Hull Lotta Damage
$ws.Range("M132").Value = -663.5792000000001  # Hull Lotta Damage
$ws.Range("H134").Value = 912.25  # Wood You Be Quiet
$ws.Range("I134").Value = 874.9  # Wood You Be Quiet
$ws.Range("K134").Value = 2624.7  # Wood You Be Quiet
$ws.Range("M134").Value = -89.69999999999982  # Wood You Be Quiet
$ws.Range("H136").Value = 1157.8  # Turali Quality
$ws.Range("I136").Value = 989.2  # Turali Quality
$ws.Range("J136").Value = 1326.4  # Turali Quality
$ws.Range("K136").Value = 2967.6  # Turali Quality
$ws.Range("L136").Value = 3979.2  # Turali Quality
$ws.Range("M136").Value = -417.6000000000004  # Turali Quality
$ws.Range("N136").Value = -9079.200000000001  # Turali Quality

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 17.333334  # Pork Is a Salty Food
$ws.Range("I2").Value = 17.333334  # Pork Is a Salty Food
$ws.Range("K2").Value = 104.000004  # Pork Is a Salty Food
$ws.Range("M2").Value = 8.999995999999996  # Pork Is a Salty Food
$ws.Range("H5").Value = 578.7857  # What a Sap
$ws.Range("I5").Value = 516.5  # What a Sap
$ws.Range("K5").Value = 1549.5  # What a Sap
$ws.Range("M5").Value = -1437.5  # What a Sap
$ws.Range("H12").Value = 99  # Butter Me Up
$ws.Range("I12").Value = 64.5  # Butter Me Up
$ws.Range("J12").Value = 110.5  # Butter Me Up
$ws.Range("K12").Value = 193.5  # Butter Me Up
$ws.Range("L12").Value = 331.5  # Butter Me Up
$ws.Range("M12").Value = -20.5  # Butter Me Up
$ws.Range("N12").Value = -677.5  # Butter Me Up
$ws.Range("H17").Value = 2980  # Chew the Fat
$ws.Range("J17").Value = 2980  # Chew the Fat
$ws.Range("L17").Value = 8940  # Chew the Fat
$ws.Range("N17").Value = -9278  # Chew the Fat
$ws.Range("H19").Value = 1136  # The Bango Zango Diet
$ws.Range("I19").Value = 50  # The Bango Zango Diet
$ws.Range("K19").Value = 150  # The Bango Zango Diet
$ws.Range("M19").Value = 24  # The Bango Zango Diet
$ws.Range("H25").Value = 1750  # Flakes for Friends
$ws.Range("I25").Value = 1000  # Flakes for Friends
$ws.Range("J25").Value = 2000  # Flakes for Friends
$ws.Range("K25").Value = 3000  # Flakes for Friends
$ws.Range("L25").Value = 6000  # Flakes for Friends
$ws.Range("M25").Value = -2831  # Flakes for Friends
$ws.Range("N25").Value = -6338  # Flakes for Friends
$ws.Range("H26").Value = 258.8  # A Grape Idea
$ws.Range("J26").Value = 248.5  # A Grape Idea
$ws.Range("L26").Value = 745.5  # A Grape Idea
$ws.Range("N26").Value = -1321.5  # A Grape Idea
$ws.Range("H30").Value = 1750  # Picnic Panic
$ws.Range("I30").Value = 1000  # Picnic Panic
$ws.Range("J30").Value = 2000  # Picnic Panic
$ws.Range("K30").Value = 3000  # Picnic Panic
$ws.Range("L30").Value = 6000  # Picnic Panic
$ws.Range("M30").Value = -2898  # Picnic Panic
$ws.Range("N30").Value = -6204  # Picnic Panic
$ws.Range("H122").Value = 947.63635  # Salt of the North
$ws.Range("J122").Value = 1990.6666  # Salt of the North
$ws.Range("L122").Value = 17915.9994  # Salt of the North
$ws.Range("N122").Value = -22815.9994  # Salt of the North
$ws.Range("H131").Value = 7825731  # The Mountain Steeped
$ws.Range("J131").Value = 14538.931  # The Mountain Steeped
$ws.Range("L131").Value = 43616.79300000001  # The Mountain Steeped
$ws.Range("N131").Value = -53696.79300000001  # The Mountain Steeped
$ws.Range("H134").Value = 1692.8096  # Dont Knock It Till Youve Tried It
$ws.Range("I134").Value = 1308.6666  # Dont Knock It Till Youve Tried It
$ws.Range("K134").Value = 3925.9998  # Dont Knock It Till Youve Tried It
$ws.Range("M134").Value = 1144.0002  # Dont Knock It Till Youve Tried It
$ws.Range("H135").Value = 578.7857  # Not-so-secret Ingredient
$ws.Range("I135").Value = 516.5  # Not-so-secret Ingredient
$ws.Range("K135").Value = 4648.5  # Not-so-secret Ingredient
$ws.Range("M135").Value = -2113.5  # Not-so-secret Ingredient
$ws.Range("H140").Value = 1854.9697  # Sweet, Sweet Bean Juice
$ws.Range("J140").Value = 2062.3704  # Sweet, Sweet Bean Juice
$ws.Range("L140").Value = 6187.111199999999  # Sweet, Sweet Bean Juice
$ws.Range("N140").Value = -16547.1112  # Sweet, Sweet Bean Juice

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7152976.5  # A Ringing Success
$ws.Range("I11").Value = 8502114  # A Ringing Success
$ws.Range("J11").Value = 70004  # A Ringing Success
$ws.Range("K11").Value = 8502114  # A Ringing Success
$ws.Range("L11").Value = 70004  # A Ringing Success
$ws.Range("M11").Value = -8501975  # A Ringing Success
$ws.Range("N11").Value = -70282  # A Ringing Success
$ws.Range("H102").Value = 1805.875  # Put the Metal to the Peddle
$ws.Range("I102").Value = 1672.8462  # Put the Metal to the Peddle
$ws.Range("J102").Value = 1963.091  # Put the Metal to the Peddle
$ws.Range("K102").Value = 1672.8462  # Put the Metal to the Peddle
$ws.Range("L102").Value = 1963.091  # Put the Metal to the Peddle
$ws.Range("M102").Value = -50.84619999999995  # Put the Metal to the Peddle
$ws.Range("N102").Value = -5207.091  # Put the Metal to the Peddle
$ws.Range("H122").Value = 2077.4211  # Awarding Academic Excellence
$ws.Range("I122").Value = 2036.3  # Awarding Academic Excellence
$ws.Range("J122").Value = 2123.111  # Awarding Academic Excellence
$ws.Range("K122").Value = 6108.9  # Awarding Academic Excellence
$ws.Range("L122").Value = 6369.333  # Awarding Academic Excellence
$ws.Range("M122").Value = -3658.9  # Awarding Academic Excellence
$ws.Range("N122").Value = -11269.333  # Awarding Academic Excellence
$ws.Range("H126").Value = 69947.13  # Gold Rush Order
$ws.Range("I126").Value = 3553.3076  # Gold Rush Order
$ws.Range("K126").Value = 10659.9228  # Gold Rush Order
$ws.Range("M126").Value = -8189.9228  # Gold Rush Order
$ws.Range("H132").Value = 3054.6758  # On Board for Lar
$ws.Range("I132").Value = 2467.1936  # On Board for Lar
$ws.Range("J132").Value = 6090  # On Board for Lar
$ws.Range("K132").Value = 7401.5808  # On Board for Lar
$ws.Range("L132").Value = 18270  # On Board for Lar
$ws.Range("M132").Value = -4871.5808  # On Board for Lar
$ws.Range("N132").Value = -23330  # On Board for Lar

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1661  # Skin off Their Backs
$ws.Range("J22").Value = 1640.125  # Skin off Their Backs
$ws.Range("L22").Value = 1640.125  # Skin off Their Backs
$ws.Range("N22").Value = -2230.125  # Skin off Their Backs
$ws.Range("H27").Value = 1661  # Fire and Hide
$ws.Range("J27").Value = 1640.125  # Fire and Hide
$ws.Range("L27").Value = 1640.125  # Fire and Hide
$ws.Range("N27").Value = -1854.125  # Fire and Hide
$ws.Range("H82").Value = 3737.9  # Trainin the Neck
$ws.Range("I82").Value = 1966.3334  # Trainin the Neck
$ws.Range("K82").Value = 1966.3334  # Trainin the Neck
$ws.Range("M82").Value = -1605.3334  # Trainin the Neck
$ws.Range("H85").Value = 3737.9  # Training Is Only Skintight (L)
$ws.Range("I85").Value = 1966.3334  # Training Is Only Skintight (L)
$ws.Range("K85").Value = 1966.3334  # Training Is Only Skintight (L)
$ws.Range("M85").Value = -718.3334  # Training Is Only Skintight (L)
$ws.Range("H122").Value = 6810.375  # Hell on Leather
$ws.Range("I122").Value = 5187  # Hell on Leather
$ws.Range("K122").Value = 15561  # Hell on Leather
$ws.Range("M122").Value = -13111  # Hell on Leather

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 16282.147  # Heavy Armoire
$ws.Range("I122").Value = 22259  # Heavy Armoire
$ws.Range("J122").Value = 1937.7  # Heavy Armoire
$ws.Range("K122").Value = 66777  # Heavy Armoire
$ws.Range("L122").Value = 5813.1  # Heavy Armoire
$ws.Range("M122").Value = -64327  # Heavy Armoire
$ws.Range("N122").Value = -10713.1  # Heavy Armoire
$ws.Range("H132").Value = 6498.5  # Comfy Cabins
$ws.Range("I132").Value = 2488  # Comfy Cabins
$ws.Range("K132").Value = 7464  # Comfy Cabins
$ws.Range("M132").Value = -4934  # Comfy Cabins
